# Updates Sheet1 with the refreshed NATMI Ptn-Ptprs TPM results.
# The header row (row 1) keeps its existing titles; rows 2-10 are rewritten
# in full, now covering every Sending/Target cluster pairing among
# ECs, FAPs and MuSCs (previously only FAPs/MuSCs senders were present).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the 16 numeric NATMI metrics (expression / specificity statistics).
$rowData = [ordered]@{
    2 = @('ECs', 'Ptn', 'Ptprs', 'ECs', 1, 0.3333333333333333, 0.06520933333333333, 0.195628, 0.007116089623850998, 0.007116089623850999, 3, 1, 2.425633666666667, 7.276901000000001, 0.0662600404061536, 0.06626004040615362, 0.1581739543142222, 1.423565588828, 0.0004715123860101775, 0.0004715123860101776)
    3 = @('ECs', 'Ptn', 'Ptprs', 'FAPs', 1, 0.3333333333333333, 0.06520933333333333, 0.195628, 0.007116089623850998, 0.007116089623850999, 3, 1, 15.503283, 46.509849, 0.4234968256437875, 0.4234968256437876, 1.010958748908, 9.098628740172, 0.003013641366697592, 0.003013641366697592)
    4 = @('ECs', 'Ptn', 'Ptprs', 'MuSCs', 1, 0.3333333333333333, 0.06520933333333333, 0.195628, 0.007116089623850998, 0.007116089623850999, 3, 1, 18.67887366666667, 56.03662100000001, 0.5102431339500588, 0.5102431339500588, 1.218036899220889, 10.962332092988, 0.003630935871143228, 0.003630935871143229)
    5 = @('FAPs', 'Ptn', 'Ptprs', 'ECs', 3, 1, 3.270036666666666, 9.81011, 0.3568488252184601, 0.3568488252184602, 3, 1, 2.425633666666667, 7.276901000000001, 0.0662600404061536, 0.06626004040615362, 7.93191102990111, 71.38719926911, 0.02364481757786361, 0.02364481757786362)
    6 = @('FAPs', 'Ptn', 'Ptprs', 'FAPs', 3, 1, 3.270036666666666, 9.81011, 0.3568488252184601, 0.3568488252184602, 3, 1, 15.503283, 46.509849, 0.4234968256437875, 0.4234968256437876, 50.69630386371, 456.26673477339, 0.1511243447147326, 0.1511243447147327)
    7 = @('FAPs', 'Ptn', 'Ptprs', 'MuSCs', 3, 1, 3.270036666666666, 9.81011, 0.3568488252184601, 0.3568488252184602, 3, 1, 18.67887366666667, 56.03662100000001, 0.5102431339500588, 0.5102431339500588, 61.08060178203446, 549.7254160383101, 0.1820796629258639, 0.1820796629258639)
    8 = @('MuSCs', 'Ptn', 'Ptprs', 'ECs', 3, 1, 5.828401, 17.485203, 0.6360350851576888, 0.6360350851576889, 3, 1, 2.425633666666667, 7.276901000000001, 0.0662600404061536, 0.06626004040615362, 14.13756568843367, 127.238091195903, 0.04214371044227981, 0.04214371044227982)
    9 = @('MuSCs', 'Ptn', 'Ptprs', 'FAPs', 3, 1, 5.828401, 17.485203, 0.6360350851576888, 0.6360350851576889, 3, 1, 15.503283, 46.509849, 0.4234968256437875, 0.4234968256437876, 90.35935014048302, 813.2341512643471, 0.2693588395623573, 0.2693588395623574)
    10 = @('MuSCs', 'Ptn', 'Ptprs', 'MuSCs', 3, 1, 5.828401, 17.485203, 0.6360350851576888, 0.6360350851576889, 3, 1, 18.67887366666667, 56.03662100000001, 0.5102431339500588, 0.5102431339500588, 108.8679659576737, 979.8116936190632, 0.3245325351530517, 0.3245325351530517)
}

foreach ($rowNum in $rowData.Keys) {
    $values = $rowData[$rowNum]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item([int]$rowNum, $col + 1).Value = $values[$col]
    }
}
